$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing row 26 (Feria Lagunitas de Puerto Montt, date 2022-04-12 / 44663,
# price 12000) needs to be preserved as a new row 27, while row 26 itself is
# updated with a newer weekly observation (date 2022-07-05 / 44747, price 14000).

# Duplicate the current row 26 into a freshly inserted row 27, preserving all
# values/styles as-is.
$ws.Rows("26").Copy()
$ws.Rows("27").Insert()

# Now overwrite row 26 in place with the new weekly price data.
$ws.Range("D26").Value = 44747
$ws.Range("K26").Value = 14000
$ws.Range("L26").Value = 14000
$ws.Range("M26").Value = 14000
$ws.Range("P26").Value = 1400
